$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 298 (shifts rows 298:326 down to 299:327)
$ws.Rows(298).Insert()

# Populate the newly inserted row 298 with the new weekly record
$ws.Cells.Item(298, 1).Value = 6
$ws.Cells.Item(298, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(298, 3).Value = "Metropolitana"
$ws.Cells.Item(298, 4).Value = 44578
$ws.Cells.Item(298, 5).Value = 13
$ws.Cells.Item(298, 6).Value = 100112032
$ws.Cells.Item(298, 7).Value = "Zapallo italiano"
$ws.Cells.Item(298, 8).Value = "Sin especificar"
$ws.Cells.Item(298, 9).Value = "Primera"
$ws.Cells.Item(298, 10).Value = 280
$ws.Cells.Item(298, 11).Value = 12000
$ws.Cells.Item(298, 12).Value = 14000
$ws.Cells.Item(298, 13).Value = 12929
$ws.Cells.Item(298, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(298, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(298, 16).Value = 259
$ws.Cells.Item(298, 17).Value = 50
$ws.Cells.Item(298, 18).Value = "Hortaliza"
